$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text content edits (sharedStrings)
# ---------------------------------------------------------------------------

# D2 (PERSON_GENDER_CODE) and E2 (PERSON_POSTCODE) keep their rich-text
# "Required" bold prefix - only the trailing sentence changes. We edit via
# Characters() so the bold run survives, then re-apply Bold/non-Bold to the
# two halves so the engine re-emits separate <r> runs.
$d2 = $ws.Range("D2")
$d2OldLen = $d2.Characters().Text.Length
$d2.Characters(9, $d2OldLen - 8).Text = ". Must be one of Not known, Male, Female, Not specified."
$d2NewLen = $d2.Characters().Text.Length
$d2.Characters(1, 8).Font.Bold = $true
$d2.Characters(9, $d2NewLen - 8).Font.Bold = $false

$e2 = $ws.Range("E2")
$e2OldLen = $e2.Characters().Text.Length
$e2.Characters(9, $e2OldLen - 8).Text = ". Must be a valid UK postcode"
$e2NewLen = $e2.Characters().Text.Length
$e2.Characters(1, 8).Font.Bold = $true
$e2.Characters(9, $e2NewLen - 8).Font.Bold = $false

# Plain-text (no rich runs) rewordings: "Required if" -> "Required only if"
$ws.Range("J2").Value = "Required only if SCHOOL_URN is 888888"
$ws.Range("O2").Value = "Optional. Must use either YYYYMMDD or DD/MM/YYYY"
$ws.Range("P2").Value = "Required only if VACCINE_GIVEN is omitted. Must be Y or N"
$ws.Range("R2").Value = "Required only if VACCINATED is N. Must be absent from school, already had elsewhere, did not attend, refused, unwell or vaccination contraindicated"
$ws.Range("V2").Value = "Required only if CARE_SETTING is 2. Must be the name of a community clinic location"

# ---------------------------------------------------------------------------
# 2. Wrap text formatting -> new cellXfs entries
#    Order matters: P2/R2/V2 (no font override) must be created before D2
#    (font override) to match xf index allocation order 6, 7.
# ---------------------------------------------------------------------------
$ws.Range("P2").WrapText = $true
$ws.Range("R2").WrapText = $true
$ws.Range("V2").WrapText = $true
$ws.Range("D2").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Row height for row 2 (wrapped text needs more vertical room)
# ---------------------------------------------------------------------------
$ws.Range("A2:Y2").RowHeight = 29.25

# ---------------------------------------------------------------------------
# 4. Column widths (best-fit-ish) for columns C:Y
# ---------------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 27.166666666666668
$ws.Columns("D").ColumnWidth = 35.736979166666664
$ws.Columns("E").ColumnWidth = 50.736979166666664
$ws.Columns("F").ColumnWidth = 52.736979166666664
$ws.Columns("G").ColumnWidth = 52.451822916666664
$ws.Columns("H").ColumnWidth = 29.877604166666668
$ws.Columns("I").ColumnWidth = 88.02213541666667
$ws.Columns("J").ColumnWidth = 34.877604166666664
$ws.Columns("K").ColumnWidth = 33.877604166666664
$ws.Columns("L").ColumnWidth = 31.736979166666668
$ws.Columns("M").ColumnWidth = 76.16666666666667
$ws.Columns("N").ColumnWidth = 14.592447916666666
$ws.Columns("O").ColumnWidth = 45.736979166666664
$ws.Columns("P").ColumnWidth = 50.736979166666664
$ws.Columns("Q").ColumnWidth = 228.02213541666666
$ws.Columns("R").ColumnWidth = 122.59244791666667
$ws.Columns("S").ColumnWidth = 7.451822916666667
$ws.Columns("T").ColumnWidth = 42.736979166666664
$ws.Columns("U").ColumnWidth = 34.022135416666664
$ws.Columns("V").ColumnWidth = 68.16666666666667
$ws.Columns("W").ColumnWidth = 32.592447916666664
$ws.Columns("X").ColumnWidth = 37.166666666666664
$ws.Columns("Y").ColumnWidth = 36.022135416666664

# ---------------------------------------------------------------------------
# 5. Sheet view: scroll so column J is left-most, select M2
# ---------------------------------------------------------------------------
[void]$ws.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = 10
